$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add round 6 results for Bálint Bakos as a new row (row 8) at the bottom
# of the Results table: Season, Round, Player, Bonus
$ws.Cells.Item(8, 1).Value = "2024"
$ws.Cells.Item(8, 2).Value = "6"
$ws.Cells.Item(8, 3).Value = "Bálint Bakos"
$ws.Cells.Item(8, 4).Value = "101"
